$d = $word.ActiveDocument

# Original (single run, trailing space kept):
#   "...parameters and adjustments. " immediately followed (same
#   paragraph) by "During the collection of the overall data, ..."
#
# Target:
#   "...parameters and adjustments." (no trailing space on that
#   sentence) + " Additionally setup the gray scaling pre-processing
#   effort in order to reduce the dimensionality of our dataset." +
#   " " (single space) + "During the collection of the overall data, ..."
#
# Net effect on the paragraph text: insert a new sentence - " Additionally
# setup the gray scaling pre-processing effort in order to reduce the
# dimensionality of our dataset." - right after "...adjustments." and
# before the pre-existing trailing space (which then naturally keeps
# separating the new sentence from "During the collection ...").

$anchor = "parameters and adjustments."

$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the sentence ending in '...adjustments.' to edit."
}

# Collapse to a zero-length point right after the period in
# "adjustments." (i.e. right before the existing trailing space) and
# insert the new sentence there.
$insertionPoint = $d.Range($findRange.End, $findRange.End)
$insertionPoint.InsertAfter(" Additionally setup the gray scaling pre-processing effort in order to reduce the dimensionality of our dataset.")

Write-Output "Edit complete."
